$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 430, shifting existing rows 430..500 down to 431..501
$ws.Rows.Item(430).Insert()

# Populate the newly inserted row 430 with the new data record
$ws.Range("A430").Value2 = 10
$ws.Range("B430").Value2 = "Vega Modelo de Temuco"
$ws.Range("C430").Value2 = "La Araucanía"
$ws.Range("D430").Value2 = 45180
$ws.Range("E430").Value2 = 9
$ws.Range("F430").Value2 = "Fruta"
$ws.Range("G430").Value2 = 100102
$ws.Range("H430").Value2 = "Cítricos"
$ws.Range("I430").Value2 = 100102006
$ws.Range("J430").Value2 = "Pomelo"
$ws.Range("K430").Value2 = "Start Ruby"
$ws.Range("L430").Value2 = "Primera"
$ws.Range("M430").Value2 = 150
$ws.Range("N430").Value2 = 15000
$ws.Range("O430").Value2 = 15000
$ws.Range("P430").Value2 = 15000
$ws.Range("Q430").Value2 = "$/bandeja 15 kilos granel"
$ws.Range("R430").Value2 = "Región de O'Higgins"
$ws.Range("S430").Value2 = 1000
$ws.Range("T430").Value2 = 15

# Ensure D430 keeps the same date number format used by the other date cells in column D
$ws.Range("D430").NumberFormat = $ws.Range("D431").NumberFormat
